$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12876.25 pesos`n✅ 12876.25 pesos = 3.32 = 945.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$cellA1 = $ws1.Range("A1")
$cellA1.Value2 = $newText

# --- Update the rate values on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value2 = 299
$ws2.Range("O10").Value2 = 3850
$ws2.Range("N12").Value2 = 3875
$ws2.Range("O12").Value2 = 284.5
